# CMPA template update: shift all campaign dates onto a single month
# (July 2023) and refresh sample group/media names, per commit:
# "aktualizace vzoru, aby datumy byly v jednom mesici"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "Groups" sheet - refresh sample group/media names + comments
#    (done first so new shared strings are interned in the same
#    order as the reference edit)
# ---------------------------------------------------------------
$groups = $wb.Worksheets.Item("Groups")

$groups.Range("A1").Value2 = "Rádio Expres, Europa 2"
$groups.Range("B1").Value2 = "Rádio Expres; Europa 2"

$commentA1 = $groups.Range("A1").Comment
[void]$commentA1.Text("Názov skupiny; musí byť rovnaký ako na záložke Spoty")

$commentB1 = $groups.Range("B1").Comment
[void]$commentB1.Text("Médiá v skupine oddelené bodkočiarkou; musia byť rovnakého mediatypu")

# ---------------------------------------------------------------
# 2) "Spots" sheet - move every date in the sample rows into July 2023
# ---------------------------------------------------------------
$spots = $wb.Worksheets.Item("Spots")

$spots.Range("F3").Value2  = 45127
$spots.Range("F4").Value2  = 45128
$spots.Range("F5").Value2  = 45129
$spots.Range("F6").Value2  = 45127
$spots.Range("F7").Value2  = 45128
$spots.Range("F8").Value2  = 45130
$spots.Range("F9").Value2  = 45126
$spots.Range("F10").Value2 = 45127
$spots.Range("F11").Value2 = 45129

$spots.Range("F12").Value2 = 45119
$spots.Range("G12").Value2 = 45127

$spots.Range("F13").Value2 = 45110
$spots.Range("G13").Value2 = 45110

$spots.Range("F14").Value2 = 45127
$spots.Range("G14").Value2 = 45129

$spots.Range("F15").Value2 = 45116
$spots.Range("G15").Value2 = 45119

$spots.Range("F16").Value2 = 45116
$spots.Range("G16").Value2 = 45119

$spots.Range("F17").Value2 = 45116
$spots.Range("G17").Value2 = 45119

$spots.Range("F18").Value2 = 45116
$spots.Range("G18").Value2 = 45119

# Sample website swapped out for a more current example
$spots.Range("B12").Value2 = "zoznam.sk"
$spots.Range("B13").Value2 = "zoznam.sk"
